$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.633.44"
$ws.Range("E2").Value = "'  +1.77%  "
$ws.Range("D3").Value = "'3.080.59"
$ws.Range("E3").Value = "'  +2.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'516.75"
$ws.Range("E5").Value = "'  +0.75%  "
$ws.Range("D6").Value = "'142.94"
$ws.Range("E6").Value = "'  +4.41%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.436"
$ws.Range("E8").Value = "'  +1.72%  "
$ws.Range("D9").Value = "'7.33"
$ws.Range("E9").Value = "'  +1.05%  "
$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "'  +1.28%  "
$ws.Range("E11").Value = "'  +2.78%  "
$ws.Range("D12").Value = "'3.599.08"
$ws.Range("E12").Value = "'  +3.02%  "
$ws.Range("E13").Value = "'  +2.67%  "
$ws.Range("D14").Value = "'25.74"
$ws.Range("E14").Value = "'  -2.42%  "
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("D16").Value = "'57.667.47"
$ws.Range("E16").Value = "'  +1.91%  "
$ws.Range("D17").Value = "'6.15"
$ws.Range("E17").Value = "'  +0.53%  "
$ws.Range("D18").Value = "'3.070.16"
$ws.Range("E18").Value = "'  +2.19%  "
$ws.Range("D19").Value = "'13.08"
$ws.Range("E19").Value = "'  +0.17%  "
$ws.Range("D20").Value = "'8.20"
$ws.Range("E20").Value = "'  +2.63%  "
$ws.Range("D21").Value = "'337.39"
$ws.Range("E21").Value = "'  +4.43%  "
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("D23").Value = "'0.501"
$ws.Range("E23").Value = "'  +0.46%  "
$ws.Range("D24").Value = "'65.54"
$ws.Range("E24").Value = "'  +2.38%  "
$ws.Range("E25").Value = "'  +5.55%  "
$ws.Range("E26").Value = "'  +0.47%  "
$ws.Range("D27").Value = "'0.0₃0934"
$ws.Range("E27").Value = "'  +6.44%  "
$ws.Range("D28").Value = "'6.47"
$ws.Range("E28").Value = "'  -1.04%  "
$ws.Range("D29").Value = "'7.08"
$ws.Range("E29").Value = "'  -1.20%  "
$ws.Range("D30").Value = "'1.81"
$ws.Range("E30").Value = "'  +0.75%  "
$ws.Range("D31").Value = "'20.85"
$ws.Range("E31").Value = "'  +1.79%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "'  -1.99%  "
$ws.Range("D33").Value = "'154.42"
$ws.Range("E33").Value = "'  +1.24%  "
$ws.Range("D34").Value = "'4.53"
$ws.Range("E34").Value = "'  +0.16%  "
$ws.Range("D35").Value = "'5.92"
$ws.Range("E35").Value = "'  +2.54%  "
$ws.Range("D36").Value = "'26.68"
$ws.Range("E36").Value = "'  +5.07%  "
$ws.Range("D37").Value = "'1.25"
$ws.Range("E37").Value = "'  +1.53%  "
$ws.Range("D38").Value = "'0.0686"
$ws.Range("E38").Value = "'  +3.45%  "
$ws.Range("D39").Value = "'3.118.92"
$ws.Range("E39").Value = "'  +2.67%  "
$ws.Range("D40").Value = "'36.97"
$ws.Range("E40").Value = "'  +0.97%  "
$ws.Range("E41").Value = "'  +1.54%  "
$ws.Range("D42").Value = "'0.672"
$ws.Range("E42").Value = "'  +3.46%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "'  -0.19%  "
$ws.Range("D44").Value = "'2.279.11"
$ws.Range("E44").Value = "'  +5.33%  "
$ws.Range("D45").Value = "'0.0254"
$ws.Range("E45").Value = "'  +5.40%  "
$ws.Range("E46").Value = "'  +1.69%  "
$ws.Range("D47").Value = "'0.958"
$ws.Range("E47").Value = "'  +1.78%  "
$ws.Range("E48").Value = "'  +4.89%  "
$ws.Range("D49").Value = "'5.87"
$ws.Range("E49").Value = "'  -3.63%  "
$ws.Range("D50").Value = "'0.0876"
$ws.Range("E50").Value = "'  +2.97%  "
$ws.Range("D51").Value = "'0.691"
$ws.Range("E51").Value = "'  +3.52%  "
